$p = $ppt.ActivePresentation

# Slide 3: "JAva vs. HAskell" -> "Java vs. Haskell"
$p.Slides.Item(3).Shapes.Item(1).TextFrame.TextRange.Text = "Java vs. Haskell"

# Slide 4: "JAva review" -> "Java review"
$p.Slides.Item(4).Shapes.Item(1).TextFrame.TextRange.Text = "Java review"

# Slide 5: "HAskell review" -> "Haskell review"
$p.Slides.Item(5).Shapes.Item(1).TextFrame.TextRange.Text = "Haskell review"

# Slide 8: "DEmo" -> "Demo"
$p.Slides.Item(8).Shapes.Item(1).TextFrame.TextRange.Text = "Demo"
